$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The document originally has 7 paragraphs:
#   1: "Took " + "Pima Indians Diabetes Database" + " from Kaggle"      (bold)
#   2: "Took " + "Indian Liver Patient Records" + " from Kaggle dataset" (bold)
#   3: "Took " + "Cleveland Heart Disease Dataset" + " from Kaggle dataset" (bold)
#   4: (empty, bold paragraph mark)
#   5: "Create conda environment to install sepereate dependecies"      (bold)
#   6: (empty, bold paragraph mark)
#   7: (empty, bold paragraph mark)
#
# Target layout:
#   1: "Step 1: Set Up the Environment"                                  (bold)
#   2: (empty, bold)
#   3: "First things first - " + "Create conda environment to install sepereate dependecies" (bold)
#   4: "Okay now when running jupyter notebook vscode will prompt which environment in the top of vscode select your environment and wait, maybe restart vscode if required"
#   5: "Then install dependencies on the environment"
#   6: (empty paragraph)
#   7: "Took " + "Pima Indians Diabetes Database" + " from Kaggle"       (bold)   [moved from old #1]
#   8: "Took " + "Indian Liver Patient Records" + " from Kaggle dataset" (bold)  [moved from old #2]
#   9: "Took " + "Cleveland Heart Disease Dataset" + " from Kaggle dataset" (bold) [moved from old #3]
#  10: (empty, bold)
#  11: "Step 2: Data Preprocessing"                                      (bold)
#  12: "Now that your data is loaded, let's preprocess it!"              (bold)
#  13: (empty, bold)
#  14: (empty, bold)   [old #4]
#  15: (empty, bold)   [old #6, untouched]
#  16: (empty, bold)   [old #7, untouched]
# ---------------------------------------------------------------------------

# Step 1: remove the three "Took ..." dataset paragraphs plus the blank
# paragraph right after them (old paragraphs 1-4). We'll recreate the three
# "Took ..." paragraphs further down, and a fresh blank paragraph near the end.
$startPos = $d.Paragraphs(1).Range.Start
$endPos = $d.Paragraphs(4).Range.End
$d.Range($startPos, $endPos).Delete()

# Paragraph 1 is now "Create conda environment to install sepereate dependecies"
# Step 2: prepend a new bold run "First things first - " to that paragraph.
$p1 = $d.Paragraphs(1)
$startOfP1 = $p1.Range
$startOfP1.Collapse(1)
$startOfP1.Font.Bold = 1
$startOfP1.InsertBefore("First things first - ")

# Step 3: after that paragraph, add the two plain (non-bold) instruction
# paragraphs.
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Okay now when running jupyter notebook vscode will prompt which environment in the top of vscode select your environment and wait, maybe restart vscode if required"
$p2.Range.Font.Bold = 0

$r = $p2.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Then install dependencies on the environment"
$p3.Range.Font.Bold = 0

# Step 4: add a blank paragraph after that.
$r = $p3.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p4 = $d.Paragraphs(4)
$p4.Range.Font.Bold = 0

# Step 5: re-add the three "Took ..." dataset paragraphs (bold), now placed
# after the blank paragraph.
$r = $p4.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p5 = $d.Paragraphs(5)
$p5.Range.Font.Bold = 1
$p5.Range.Text = "Took Pima Indians Diabetes Database from Kaggle"

$r = $p5.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p6 = $d.Paragraphs(6)
$p6.Range.Font.Bold = 1
$p6.Range.Text = "Took Indian Liver Patient Records from Kaggle dataset"

$r = $p6.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p7 = $d.Paragraphs(7)
$p7.Range.Font.Bold = 1
$p7.Range.Text = "Took Cleveland Heart Disease Dataset from Kaggle dataset"

# Step 6: empty bold paragraph, then "Step 2: Data Preprocessing", then the
# "Now that your data is loaded..." paragraph, then one more blank bold
# paragraph.
$r = $p7.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p8 = $d.Paragraphs(8)
$p8.Range.Font.Bold = 1

$r = $p8.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p9 = $d.Paragraphs(9)
$p9.Range.Font.Bold = 1
$p9.Range.Text = "Step 2: Data Preprocessing"

$r = $p9.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p10 = $d.Paragraphs(10)
$p10.Range.Font.Bold = 1
$p10.Range.Text = "Now that your data is loaded, let" + [char]0x2019 + "s preprocess it!"

$r = $p10.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p11 = $d.Paragraphs(11)
$p11.Range.Font.Bold = 1

# Step 7: insert the title paragraph ("Step 1: Set Up the Environment")
# and a blank bold paragraph before everything we just built.
$first = $d.Paragraphs(1)
$startRng = $first.Range
$startRng.Collapse(1)
$startRng.InsertParagraphBefore() | Out-Null
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Font.Bold = 1
$titlePara.Range.Text = "Step 1: Set Up the Environment"

$r = $titlePara.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$blankPara = $d.Paragraphs(2)
$blankPara.Range.Font.Bold = 1

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
foreach ($p in $d.Paragraphs) {
    Write-Host "[$($p.Range.Text)]"
}
